# Update the Metadata sheet of the wh-payer-medicare-coverage CodeSystem
# for the new IG build (Version 6.0.0, published 2022-01-21):
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to the new publication timestamp
#  - Publisher gets a value ("Alvearie Team")
#  - the duplicated "Contact" row is dropped and replaced with "Jurisdiction"
#  - Case Sensitive gets a value ("true")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version
$ws.Range("B3").Value = "6.0.0"

# Date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws.Range("B9").Value = "Alvearie Team"

# Row 11 duplicated row 10 (Contact / No display for ContactDetail) - remove it,
# shifting the rows below up by one.
$ws.Rows.Item(11).Delete()

# Row 10 (was "Contact" / "No display for ContactDetail") becomes "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive now has a value. Plain assignment of the literal text "true"
# gets auto-coerced to a Boolean by Excel's usual TRUE/FALSE parsing, so
# instead stage the text in a scratch cell (via a formula that yields the
# string "true"), copy it, and paste-special just the value into B14 - this
# keeps the cell's text type and its existing style untouched.
$scratch = $ws.Range("F1")
$scratch.Formula = "=""true"""
$scratch.Copy()
$caseSensitiveCell = $ws.Range("B14")
$caseSensitiveCell.PasteSpecial(-4163, -4142, $false, $false)
$scratch.Clear()
